# Flip the sign (positive -> negative) of the "material recycled" values
# for each metal/component on every yearly sheet from 2002 through 2100.
# The affected cells on each sheet are:
#   C2 (Nd / Generator Offshore)
#   B4 (Cu / Generator Onshore)
#   C4 (Cu / Generator Offshore)
#   E4 (Cu / Wires)
#   D5 (Si / Panel)
# Sheets "2000" and "2001" only contain zeros in these cells, so they are
# left untouched (negating zero is a no-op and the original diff does not
# touch them).

$wb = $excel.ActiveWorkbook

$targetCells = @("C2", "B4", "C4", "E4", "D5")

for ($year = 2002; $year -le 2100; $year++) {
    $ws = $wb.Worksheets.Item("$year")

    foreach ($addr in $targetCells) {
        $cell = $ws.Range($addr)
        $current = $cell.Value2
        if ($current -ne 0) {
            $cell.Value = -$current
        }
    }
}
